$wb = $excel.ActiveWorkbook

$wsYearly = $wb.Worksheets.Item("Yearly")
$wsAllTime = $wb.Worksheets.Item("All Time")

# Update the taxable dividend value for 2017 which cascades through the
# dependent SUM/formula cells on both sheets (totals row and grand total).
$wsYearly.Range("L3").Value = 57.62

# Restore the view/selection state on each sheet.
$wsYearly.Activate()
$wsYearly.Range("F37").Select()

$wsAllTime.Activate()
$wsAllTime.Application.ActiveWindow.ScrollRow = 1
$wsAllTime.Range("N37").Select()
